$d = $word.ActiveDocument

# Pairs of (old, new) division-fact strings to update. Order matters: a
# couple of the new values coincide with old values used elsewhere in the
# document, so each replacement is executed once and in the same
# top-to-bottom order the cells appear in the document, which guarantees
# every "old" substring is rewritten before it could be (re)created by a
# later replacement.
$pairs = @(
    @("79÷7=", "11÷6="),
    @("56÷8=", "64÷7="),
    @("53÷6=", "55÷3="),
    @("11÷3=", "46÷2="),
    @("75÷5=", "59÷7="),
    @("69÷4=", "35÷6="),
    @("56÷6=", "12÷6="),
    @("94÷8=", "55÷7="),
    @("21÷5=", "88÷2="),
    @("84÷2=", "69÷4="),
    @("74÷9=", "56÷4="),
    @("12÷5=", "65÷5="),
    @("69÷3=", "91÷6="),
    @("97÷6=", "46÷9="),
    @("38÷6=", "97÷8="),
    @("42÷9=", "95÷5="),
    @("57÷9=", "51÷8="),
    @("94÷4=", "25÷5="),
    @("85÷8=", "74÷7="),
    @("37÷4=", "40÷6="),
    @("16÷7=", "93÷8="),
    @("27÷3=", "26÷7="),
    @("56÷7=", "45÷8="),
    @("36÷2=", "31÷5="),
    @("73÷5=", "20÷7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
